$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for every existing data row
# (rows 2-452) from 2023-09-06 (45175) to 2023-09-08 (45177).
$ws.Range("C2:C452").Value = 45177

# Row 452 gains an explicit row height (15, customHeight) in the target file.
$ws.Rows.Item(452).RowHeight = 15

# Append the new record as row 453.
$ws.Range("A453").Value = "A 41844-2023"

$ws.Range("B453").Value = 45176
$ws.Range("B453").NumberFormat = "YYYY-MM-DD"

$ws.Range("C453").Value = 45177
$ws.Range("C453").NumberFormat = "YYYY-MM-DD"

$ws.Range("D453").Value = "GÄVLEBORGS LÄN"
$ws.Range("E453").Value = "GÄVLE"

$ws.Range("G453").Value = 0.9
$ws.Range("H453").Value = 0
$ws.Range("I453").Value = 0
$ws.Range("J453").Value = 0
$ws.Range("K453").Value = 0
$ws.Range("L453").Value = 0
$ws.Range("M453").Value = 0
$ws.Range("N453").Value = 0
$ws.Range("O453").Value = 0
$ws.Range("P453").Value = 0
$ws.Range("Q453").Value = 0

# R column keeps the same empty, wrap-text styled cell used throughout the sheet.
$ws.Range("R453").WrapText = $true
